$d = $word.ActiveDocument

# Right single quotation mark (U+2019) used in "haven't" below.
$rsquo = [char]0x2019

# Delete whole paragraphs (text + trailing paragraph mark) that were
# removed in the commit "More work with BFS, DFS".
$toRemove = @(
    "Word Break",
    "Coin Change",
    "Reverse",
    ("2Sum (lol, haven" + $rsquo + "t done this yet for some reason)")
)

foreach ($text in $toRemove) {
    $find = $d.Content.Find
    $find.ClearFormatting()
    $result = $find.Execute($text + "^p", $false, $false, $false, $false, $false, $true, 1, $false, "", 2)
    if (-not $result) {
        Write-Output "NOT FOUND: $text"
    }
}
